$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the previously-empty header row (row 1) with the API field labels,
# matching the formatting used for column B (bold Segoe UI 7, white fill)
# but without the thin border used by the data rows.
$ws.Range("B2").Copy()
$ws.Range("C1:P1").PasteSpecial(-4122)
$ws.Range("C1:P1").Borders.LineStyle = -4142
$excel.CutCopyMode = $false

$headers = @(
    "API field",
    "flow name",
    "UUID dataset",
    "flow property",
    "unit",
    "type",
    "is reference flow?",
    "in/out",
    "sugarcane cultivation and transport to refinery",
    "bioethanol production, from sugarcane",
    "soybean cultivation and transport to refinery",
    "biodiesel production, from soybean",
    "rape seed cultivation and transport to refinery",
    "biodiesel production, from rape seed"
)

$col = 3
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$ws.Rows.Item(1).RowHeight = 48

$null = $ws.Range("S3").Select()
